$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E6").Value = 49

$ws.Range("E8").Value = 37
$ws.Range("F8").Value = 13
$ws.Range("H8").Value = 13

$ws.Range("E10").Value = 20

$ws.Range("E16").Value = 282

$ws.Range("E18").Value = 84
